$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Mark rows 10, 11, 12, 14 (column A) as FINALIZADO, copying the
#    "FINALIZADO" formatting (green fill) from A2.
$src = $ws.Range("A2")
$targets = @(10, 11, 12, 14)
foreach ($r in $targets) {
    $cell = $ws.Cells.Item($r, 1)
    $src.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cell.Value = "FINALIZADO"
}

# 2) Remove the "Migrar pro Vue.js" task row entirely (old row 18),
#    shifting everything below it up by one row.
$ws.Rows.Item(18).Delete() | Out-Null

# 3) Leave selection on C17 to match the saved view state.
$ws.Range("C17").Select() | Out-Null
